# check #11 on 31/3/2025
# Remove villages that no longer belong in the "Kitagwenda / 2025 / A" block
# and simplify the missing_trainings text for the remaining rows that were
# affected by that cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the rows for villages that were dropped -------------------
# (village -> original row number, deleted highest-first so row numbers of
# rows still to be removed don't shift while we work)
$rowsToDelete = @(22, 21, 18, 13, 12, 6, 5, 4)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# --- 2. Simplify the "missing_trainings" text for the surviving rows -----
# Kitagwenda / 2025 / A rows that used to combine "HOR 1" with the coffee
# training now only show the coffee training.
$simplifiedCoffeeRows = @(4, 9, 12, 13)
foreach ($r in $simplifiedCoffeeRows) {
    $ws.Cells.Item($r, 5).Value = "Coffee Champions - ToT 1"
}

# Rakai / 2025 / A rows (Kammengo .. Nnongo_A) all collapse down to the
# single HHT 2 training instead of the various combined lists.
for ($r = 20; $r -le 36; $r++) {
    $ws.Cells.Item($r, 5).Value = "HHT 2 - Financial Literacy and VSLA and and Business Enterprise Selection"
}
